$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style of the existing header row (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Cells.Item(1,8).Value = "Save"

# New column H data rows 2-9, all zeros (numeric, unstyled like the rest of the data)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
